# Update "想去人数" (want-to-attend count) figures in the 展览 (F column)
# and propagate the same updates to the merged 全部类型 sheet, matching
# the data refresh captured at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll     = $wb.Worksheets.Item("全部类型")

# Each entry: row in 展览 sheet, row in 全部类型 sheet, new value
$updates = @(
    @{ RowA = 3;  RowB = 5;  New = 852 },
    @{ RowA = 4;  RowB = 6;  New = 258 },
    @{ RowA = 5;  RowB = 8;  New = 75 },
    @{ RowA = 7;  RowB = 10; New = 776 },
    @{ RowA = 9;  RowB = 12; New = 1488 },
    @{ RowA = 11; RowB = 14; New = 1047 },
    @{ RowA = 13; RowB = 16; New = 70 },
    @{ RowA = 14; RowB = 17; New = 197 },
    @{ RowA = 15; RowB = 18; New = 55 },
    @{ RowA = 17; RowB = 20; New = 48 },
    @{ RowA = 18; RowB = 21; New = 37 },
    @{ RowA = 22; RowB = 30; New = 567 },
    @{ RowA = 24; RowB = 32; New = 36 },
    @{ RowA = 26; RowB = 34; New = 770 },
    @{ RowA = 27; RowB = 35; New = 255 }
)

foreach ($u in $updates) {
    $wsExhibit.Range("F" + $u.RowA).Value = $u.New
    $wsAll.Range("F" + $u.RowB).Value = $u.New
}
